$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(71).Insert()

$ws.Range("A71").Value = 3
$ws.Range("B71").Value = "Femacal de La Calera"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44629
$ws.Range("E71").Value = 5
$ws.Range("F71").Value = 100112052
$ws.Range("G71").Value = "Albahaca"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 115
$ws.Range("K71").Value = 5000
$ws.Range("L71").Value = 5500
$ws.Range("M71").Value = 5239
$ws.Range("N71").Value = "$/docena de matas"
$ws.Range("O71").Value = "Provincia de Quillota"
$ws.Range("P71").Value = 873
$ws.Range("Q71").Value = 6
$ws.Range("R71").Value = "Hortaliza"
